$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Video sizes are now computed instead of defined in the layout xml file:
# add two more rows of pre-computed "small" profile values below the
# existing shared-formula table (A3:G11), mirroring rows 4/5 ("small").

# Row 13 -> corresponds to a 426x160 "small" screen (like row 4, but at 160dpi)
$ws.Range("A13").Value = "small"
$ws.Range("B13").Value = 426
$ws.Range("C13").Value = 160
$ws.Range("D13").Value = 0.75
$ws.Range("E13").Value = 319.5
$ws.Range("F13").Value = 2.6625000000000001
$ws.Range("G13").Value = 3.330000938

# Row 14 -> corresponds to a 320x160 "small" screen (like row 5, but at 160dpi)
$ws.Range("A14").Value = "small"
$ws.Range("B14").Value = 320
$ws.Range("C14").Value = 160
$ws.Range("D14").Value = 0.75
$ws.Range("E14").Value = 240
$ws.Range("F14").Value = 2

[void]$ws.Range("E5").Select()

Write-Output ("Dimension: " + $ws.UsedRange.Address())
